# Update "想去人数" (number of people interested) figures across sheets
# following the upstream data refresh (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 271
$ws1.Range("F7").Value = 4938
$ws1.Range("F8").Value = 4938
$ws1.Range("F9").Value = 48
$ws1.Range("F11").Value = 481
$ws1.Range("F14").Value = 670
$ws1.Range("F15").Value = 4673
$ws1.Range("F17").Value = 195
$ws1.Range("F18").Value = 88
$ws1.Range("F19").Value = 232
$ws1.Range("F20").Value = 3645
$ws1.Range("F23").Value = 30
$ws1.Range("F24").Value = 3398
$ws1.Range("F29").Value = 219
$ws1.Range("F32").Value = 89
$ws1.Range("F34").Value = 54
$ws1.Range("F36").Value = 5973
$ws1.Range("F37").Value = 939
$ws1.Range("F38").Value = 451
$ws1.Range("F42").Value = 1234
$ws1.Range("F43").Value = 131
$ws1.Range("F44").Value = 571
$ws1.Range("F46").Value = 2096
$ws1.Range("F49").Value = 737

# -----------------------------------------------------------------
# Sheet "演出" (Performances)
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 46
$ws2.Range("F9").Value = 70
$ws2.Range("F20").Value = 43
$ws2.Range("F23").Value = 780

# -----------------------------------------------------------------
# Sheet "全部类型" (All types - combined listing)
# -----------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 271
$ws4.Range("F9").Value = 4938
$ws4.Range("F10").Value = 4938
$ws4.Range("F11").Value = 48
$ws4.Range("F12").Value = 46
$ws4.Range("F14").Value = 70
$ws4.Range("F15").Value = 481
$ws4.Range("F17").Value = 670
$ws4.Range("F18").Value = 4673
$ws4.Range("F20").Value = 195
$ws4.Range("F21").Value = 88
$ws4.Range("F22").Value = 232
$ws4.Range("F23").Value = 3645
$ws4.Range("F24").Value = 3398
$ws4.Range("F28").Value = 219
$ws4.Range("F31").Value = 89
$ws4.Range("F32").Value = 54
$ws4.Range("F35").Value = 5973
$ws4.Range("F36").Value = 939
$ws4.Range("F37").Value = 451
$ws4.Range("F42").Value = 1234
$ws4.Range("F43").Value = 131
$ws4.Range("F44").Value = 571
$ws4.Range("F45").Value = 2096
$ws4.Range("F48").Value = 737
